# Add two new rows (Anapa & Gelendzhik) as new local extremums, matching
# the formatting style of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- formatting (apply first so no brand-new styles get synthesized) -------
foreach ($r in 9, 10) {
    $fullRow = $ws.Range("A$r" + ":T$r")
    $fullRow.HorizontalAlignment = -4108   # xlCenter, matches existing rows

    $ws.Range("F$r").NumberFormat = "0.000"
    $ws.Range("G$r").NumberFormat = "0"
    $ws.Range("H$r" + ":T$r").NumberFormat = "0.000"
}

# --- row 9: город-курорт Анапа ---------------------------------------------
$ws.Cells.Item(9, 1).Value = [double]"3703000"
$ws.Cells.Item(9, 2).Value = "город-курорт Анапа"
$ws.Cells.Item(9, 3).Value = [double]"2021"
$ws.Cells.Item(9, 4).Value = [double]"1763"
$ws.Cells.Item(9, 5).Value = [double]"212839"
$ws.Cells.Item(9, 6).Value = [double]"0.13402618880938175"
$ws.Cells.Item(9, 7).Value = [double]"23354.08466"
$ws.Cells.Item(9, 8).Value = [double]"0.66409633572794458"
$ws.Cells.Item(9, 9).Value = [double]"0.11628507933226523"
$ws.Cells.Item(9, 10).Value = [double]"101.42730103688706"
$ws.Cells.Item(9, 11).Value = [double]"38.700000000000003"
$ws.Cells.Item(9, 12).Value = [double]"2.2129402975958352E-3"
$ws.Cells.Item(9, 13).Value = [double]"1.3343419204187201E-3"
$ws.Cells.Item(9, 14).Value = [double]"4.5414609164673766E-3"
$ws.Cells.Item(9, 15).Value = [double]"0.87642302397586913"
$ws.Cells.Item(9, 16).Value = [double]"0.41760062770450901"
$ws.Cells.Item(9, 17).Value = [double]"16.504679212456363"
$ws.Cells.Item(9, 18).Value = [double]"1.9733225583657132E-4"
$ws.Cells.Item(9, 19).Value = [double]"4.0608159218940135E-2"
$ws.Cells.Item(9, 20).Value = [double]"196.87771655772673"

# --- row 10: город-курорт Геленджик -----------------------------------------
$ws.Cells.Item(10, 1).Value = [double]"3708000"
$ws.Cells.Item(10, 2).Value = "город-курорт Геленджик"
$ws.Cells.Item(10, 3).Value = [double]"2021"
$ws.Cells.Item(10, 4).Value = [double]"-43"
$ws.Cells.Item(10, 5).Value = [double]"115048"
$ws.Cells.Item(10, 6).Value = [double]"0.16380119602252974"
$ws.Cells.Item(10, 7).Value = [double]"25688.450580000001"
$ws.Cells.Item(10, 8).Value = [double]"1.2728165635213129"
$ws.Cells.Item(10, 9).Value = [double]"0.25154718030734996"
$ws.Cells.Item(10, 10).Value = [double]"123.63113805385579"
$ws.Cells.Item(10, 11).Value = [double]"42.7"
$ws.Cells.Item(10, 12).Value = [double]"2.6684514289687781E-3"
$ws.Cells.Item(10, 13).Value = [double]"5.3369028579375562E-3"
$ws.Cells.Item(10, 14).Value = [double]"3.4689868576594121E-3"
$ws.Cells.Item(10, 15).Value = [double]"0.11611675126903553"
$ws.Cells.Item(10, 16).Value = [double]"5.9559488213615182E-2"
$ws.Cells.Item(10, 17).Value = [double]"3.6999218004658925"
$ws.Cells.Item(10, 18).Value = [double]"2.1730060496488423E-4"
$ws.Cells.Item(10, 19).Value = [double]"5.0326820109867186E-2"
$ws.Cells.Item(10, 20).Value = [double]"107.78377544607469"

# --- move / refresh the active selection like the author left it -----------
$ws.Range("B16").Select() | Out-Null
